$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column numeric-looking values stay formatted as text (preserve trailing zeros / multi-dot format)
$textCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D18","D21","D22","D23","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D43","D44","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.530.73'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.481.66'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '313.40'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = '92.43'
$ws.Range("E6").Value = '  -2.65%  '
$ws.Range("D7").Value = '0.547'
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +2.45%  '
$ws.Range("D10").Value = '32.88'
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("D11").Value = '0.0792'
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '2.862.93'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '16.42'
$ws.Range("E14").Value = '  +9.50%  '
$ws.Range("D15").Value = '6.92'
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").Value = '2.489.87'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '0.775'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").Value = '41.549.32'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("E19").Value = '  +2.98%  '
$ws.Range("E20").Value = '  +2.40%  '
$ws.Range("D21").Value = '72.38'
$ws.Range("E21").Value = '  +5.36%  '
$ws.Range("D22").Value = '11.21'
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("D23").Value = '236.84'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '24.81'
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("D30").Value = '35.94'
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("D31").Value = '157.86'
$ws.Range("E31").Value = '  +3.73%  '
$ws.Range("D32").Value = '5.48'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").Value = '  -0.97%  '
$ws.Range("D34").Value = '0.0757'
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("D35").Value = '17.42'
$ws.Range("E35").Value = '  +1.70%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.106'
$ws.Range("E36").Value = '  +3.86%  '
$ws.Range("B37").Value = 'ApeXProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D37").Value = '2.37'
$ws.Range("E37").Value = '  -10.85%  '
$ws.Range("E38").Value = '  -4.75%  '
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("E41").Value = '  -4.44%  '
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").Value = '1.974.03'
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("D44").Value = '19.23'
$ws.Range("E44").Value = '  -4.08%  '
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("D47").Value = '8.96'
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("D48").Value = '2.720.20'
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").Value = '98.12'
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("D50").Value = '68.30'
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("D51").Value = '72.52'
$ws.Range("E51").Value = '  -3.57%  '
